$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: new CNPJ value, and the number format reverts to General (while font
# stays the same grey Arial look it already had).
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Color = 5656909
$ws.Range("A2").Value = 29148959000150

# New row 6: another CNPJ value, using the same default numeric style as
# rows like A5 (integer format, default font -> s="1").
$ws.Range("A6").NumberFormat = "0"
$ws.Range("A6").Value = 51041667000173

# Selection moves to A2.
[void]$ws.Range("A2").Select()
